# Edit DaniGargya_MA_proposal_Mar24.docx per the commit diff:
#  1. "...started the project in January 2023, finished the main project
#     phase in July 23. Both are secondary schools..." - the trailing
#     run break right after "July 23." collapses (text unchanged).
#  2. "Emmendingen, a suburb of Freiburg." -> "Emmendingen, a neighbouring
#     city of Freiburg."
#  3. "is a Montessori center" -> "is a Montessori centre" (British
#     spelling), dropping the now-stale spell-check proofErr markers
#     that wrapped "center".

$d = $word.ActiveDocument

# 1) + 2) "a suburb of Freiburg" -> "a neighbouring city of Freiburg" and
#    "Montessori center " -> "Montessori centre " (trailing space pulls
#    the orphaned spellStart/spellEnd proofErr pair for "center" into the
#    replaced range so they get dropped along with the stale spelling).
$d.Content.Find.Execute(
    ", a suburb of Freiburg. The second school, ANGELL Schule, is a Montessori center ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", a neighbouring city of Freiburg. The second school, ANGELL Schule, is a Montessori centre ",
    2) | Out-Null
